# Applies the "Automatic update of files" edit:
#  1. Column C ("Förändrad" / last-changed date) on every data row moves
#     from serial date 45184 (2023-09-15) to 45186 (2023-09-17).
#  2. Every existing single-argument HYPERLINK(...) formula (columns
#     S, T, V, W, X, Y on the rows that have species-find links) gains a
#     second argument: the row's case id (column A), used as the
#     link's friendly display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count - 1   # UsedRange starts at row 0 in this sheet

# --- 1. Bump the "Förändrad" date in column C for every data row -----------
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# --- 2. Add the friendly-name argument to the existing HYPERLINK formulas --
$linkCols = @("S", "T", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    $caseId = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $linkCols) {
        $addr = "$col$r"
        $cell = $ws.Range($addr)
        $formula = $cell.Formula

        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        if ($formula.ToUpper().IndexOf("HYPERLINK(") -lt 0) {
            continue
        }
        # Already has a second (display-text) argument -> nothing to do
        # (single-argument HYPERLINK("url") never contains a comma).
        if ($formula.IndexOf(",") -ge 0) {
            continue
        }

        # Only touch formulas of the exact shape HYPERLINK("url") -> add
        # the display-text argument, turning it into HYPERLINK("url", "id").
        if ($formula.TrimEnd().EndsWith('")')) {
            $newFormula = $formula.TrimEnd()
            $newFormula = $newFormula.Substring(0, $newFormula.Length - 1) + ', "' + $caseId + '")'
            $cell.Formula = $newFormula
        }
    }
}
